{"js": "// Update the two-digit multiplication table: each cell's text is\n// replaced with the new problem text per the diff (values given in\n// document/row-major order, matching table.values' shape).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst newValues = [\n  [\"10\u00d797=\", \"78\u00d720=\", \"64\u00d730=\", \"60\u00d783=\", \"84\u00d764=\"],\n  [\"90\u00d799=\", \"87\u00d791=\", \"40\u00d743=\", \"35\u00d717=\", \"92\u00d776=\"],\n  [\"91\u00d784=\", \"28\u00d718=\", \"29\u00d719=\", \"15\u00d798=\", \"84\u00d762=\"],\n  [\"69\u00d717=\", \"23\u00d793=\", \"99\u00d757=\", \"14\u00d781=\", \"51\u00d757=\"],\n  [\"90\u00d774=\", \"58\u00d797=\", \"15\u00d751=\", \"69\u00d781=\", \"57\u00d767=\"],\n  [\"95\u00d752=\", \"24\u00d751=\", \"61\u00d750=\", \"80\u00d763=\", \"53\u00d734=\"],\n  [\"47\u00d722=\", \"30\u00d736=\", \"53\u00d797=\", \"10\u00d716=\", \"86\u00d710=\"],\n  [\"99\u00d759=\", \"65\u00d770=\", \"25\u00d796=\", \"16\u00d767=\", \"14\u00d779=\"],\n  [\"32\u00d787=\", \"14\u00d751=\", \"92\u00d790=\", \"10\u00d722=\", \"14\u00d735=\"],\n  [\"65\u00d756=\", \"38\u00d765=\", \"75\u00d781=\", \"29\u00d757=\", \"99\u00d797=\"],\n  [\"53\u00d722=\", \"41\u00d717=\", \"95\u00d788=\", \"10\u00d720=\", \"39\u00d754=\"],\n  [\"77\u00d715=\", \"95\u00d763=\", \"61\u00d729=\", \"14\u00d734=\", \"60\u00d796=\"],\n  [\"55\u00d758=\", \"39\u00d732=\", \"35\u00d731=\", \"90\u00d766=\", \"31\u00d720=\"],\n  [\"13\u00d725=\", \"23\u00d744=\", \"23\u00d769=\", \"43\u00d785=\", \"22\u00d788=\"],\n  [\"67\u00d739=\", \"74\u00d784=\", \"86\u00d735=\", \"99\u00d793=\", \"52\u00d761=\"],\n  [\"27\u00d718=\", \"65\u00d798=\", \"85\u00d758=\", \"92\u00d754=\", \"99\u00d719=\"],\n  [\"65\u00d783=\", \"88\u00d716=\", \"49\u00d765=\", \"93\u00d754=\", \"23\u00d756=\"],\n  [\"35\u00d783=\", \"47\u00d788=\", \"57\u00d757=\", \"58\u00d740=\", \"87\u00d777=\"],\n  [\"21\u00d758=\", \"42\u00d766=\", \"39\u00d718=\", \"73\u00d793=\", \"98\u00d763=\"],\n  [\"28\u00d741=\", \"22\u00d723=\", \"58\u00d753=\", \"28\u00d765=\", \"64\u00d728=\"]\n];\n\nif (table.rowCount !== newValues.length) {\n  throw new Error(\"Unexpected row count: \" + table.rowCount);\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the two-digit multiplication table: each cell's text is\n# replaced with the new problem text per the diff (values given in\n# document/row-major order matching Table.Cell(row, col) addressing).\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"10\u00d797=\", \"78\u00d720=\", \"64\u00d730=\", \"60\u00d783=\", \"84\u00d764=\"),\n    @(\"90\u00d799=\", \"87\u00d791=\", \"40\u00d743=\", \"35\u00d717=\", \"92\u00d776=\"),\n    @(\"91\u00d784=\", \"28\u00d718=\", \"29\u00d719=\", \"15\u00d798=\", \"84\u00d762=\"),\n    @(\"69\u00d717=\", \"23\u00d793=\", \"99\u00d757=\", \"14\u00d781=\", \"51\u00d757=\"),\n    @(\"90\u00d774=\", \"58\u00d797=\", \"15\u00d751=\", \"69\u00d781=\", \"57\u00d767=\"),\n    @(\"95\u00d752=\", \"24\u00d751=\", \"61\u00d750=\", \"80\u00d763=\", \"53\u00d734=\"),\n    @(\"47\u00d722=\", \"30\u00d736=\", \"53\u00d797=\", \"10\u00d716=\", \"86\u00d710=\"),\n    @(\"99\u00d759=\", \"65\u00d770=\", \"25\u00d796=\", \"16\u00d767=\", \"14\u00d779=\"),\n    @(\"32\u00d787=\", \"14\u00d751=\", \"92\u00d790=\", \"10\u00d722=\", \"14\u00d735=\"),\n    @(\"65\u00d756=\", \"38\u00d765=\", \"75\u00d781=\", \"29\u00d757=\", \"99\u00d797=\"),\n    @(\"53\u00d722=\", \"41\u00d717=\", \"95\u00d788=\", \"10\u00d720=\", \"39\u00d754=\"),\n    @(\"77\u00d715=\", \"95\u00d763=\", \"61\u00d729=\", \"14\u00d734=\", \"60\u00d796=\"),\n    @(\"55\u00d758=\", \"39\u00d732=\", \"35\u00d731=\", \"90\u00d766=\", \"31\u00d720=\"),\n    @(\"13\u00d725=\", \"23\u00d744=\", \"23\u00d769=\", \"43\u00d785=\", \"22\u00d788=\"),\n    @(\"67\u00d739=\", \"74\u00d784=\", \"86\u00d735=\", \"99\u00d793=\", \"52\u00d761=\"),\n    @(\"27\u00d718=\", \"65\u00d798=\", \"85\u00d758=\", \"92\u00d754=\", \"99\u00d719=\"),\n    @(\"65\u00d783=\", \"88\u00d716=\", \"49\u00d765=\", \"93\u00d754=\", \"23\u00d756=\"),\n    @(\"35\u00d783=\", \"47\u00d788=\", \"57\u00d757=\", \"58\u00d740=\", \"87\u00d777=\"),\n    @(\"21\u00d758=\", \"42\u00d766=\", \"39\u00d718=\", \"73\u00d793=\", \"98\u00d763=\"),\n    @(\"28\u00d741=\", \"22\u00d723=\", \"58\u00d753=\", \"28\u00d765=\", \"64\u00d728=\")\n)\n\n$rowCount = $newValues.Count\nif ($tbl.Rows.Count -ne $rowCount) {\n    throw \"Unexpected row count: $($tbl.Rows.Count)\"\n}\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $tbl.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
